$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 436.875
$ws.Range("I2").Value = 249.16667
$ws.Range("K2").Value = 249.16667
$ws.Range("M2").Value = -136.16667
$ws.Range("H38").Value = 255
$ws.Range("J38").Value = 500
$ws.Range("L38").Value = 1500
$ws.Range("N38").Value = -2244
$ws.Range("H39").Value = 155
$ws.Range("I39").Value = 107.57143
$ws.Range("K39").Value = 322.71429
$ws.Range("M39").Value = -26.71429000000001
$ws.Range("H41").Value = 682.9167
$ws.Range("I41").Value = 662.25
$ws.Range("J41").Value = 724.25
$ws.Range("K41").Value = 662.25
$ws.Range("L41").Value = 724.25
$ws.Range("M41").Value = -222.25
$ws.Range("N41").Value = -1604.25
$ws.Range("H76").Value = 4978.25
$ws.Range("I76").Value = 4003
$ws.Range("J76").Value = 6170.222
$ws.Range("K76").Value = 4003
$ws.Range("L76").Value = 6170.222
$ws.Range("M76").Value = -3688
$ws.Range("N76").Value = -6800.222
$ws.Range("H79").Value = 4978.25
$ws.Range("I79").Value = 4003
$ws.Range("J79").Value = 6170.222
$ws.Range("K79").Value = 4003
$ws.Range("L79").Value = 6170.222
$ws.Range("M79").Value = -2911
$ws.Range("N79").Value = -8354.222
$ws.Range("H80").Value = 2990.3784
$ws.Range("I80").Value = 2269
$ws.Range("J80").Value = 3839.0588
$ws.Range("K80").Value = 6807
$ws.Range("L80").Value = 11517.1764
$ws.Range("M80").Value = -5809
$ws.Range("N80").Value = -13513.1764
$ws.Range("H83").Value = 2990.3784
$ws.Range("I83").Value = 2269
$ws.Range("J83").Value = 3839.0588
$ws.Range("K83").Value = 20421
$ws.Range("L83").Value = 34551.5292
$ws.Range("M83").Value = -15429
$ws.Range("N83").Value = -44535.5292
$ws.Range("H88").Value = 5065.0835
$ws.Range("I88").Value = 5571.2856
$ws.Range("J88").Value = 4356.4
$ws.Range("K88").Value = 5571.2856
$ws.Range("L88").Value = 4356.4
$ws.Range("M88").Value = -5165.2856
$ws.Range("N88").Value = -5168.4
$ws.Range("H91").Value = 5065.0835
$ws.Range("I91").Value = 5571.2856
$ws.Range("J91").Value = 4356.4
$ws.Range("K91").Value = 5571.2856
$ws.Range("L91").Value = 4356.4
$ws.Range("M91").Value = -4167.2856
$ws.Range("N91").Value = -7164.4
$ws.Range("H92").Value = 656.52
$ws.Range("I92").Value = 747.5294
$ws.Range("K92").Value = 747.5294
$ws.Range("M92").Value = 500.4706
$ws.Range("H98").Value = 1185.9286
$ws.Range("I98").Value = 1215.037
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 1215.037
$ws.Range("L98").Value = 400
$ws.Range("M98").Value = 282.963
$ws.Range("N98").Value = -3396
$ws.Range("H122").Value = 1185.9286
$ws.Range("I122").Value = 1215.037
$ws.Range("J122").Value = 400
$ws.Range("K122").Value = 3645.111
$ws.Range("L122").Value = 1200
$ws.Range("M122").Value = -1195.111
$ws.Range("N122").Value = -6100
$ws.Range("H132").Value = 9423.462
$ws.Range("I132").Value = 9892.083000000001
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 29676.249
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -27146.249
$ws.Range("N132").Value = -16460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3899.2
$ws.Range("I45").Value = 4128.4287
$ws.Range("J45").Value = 3698.625
$ws.Range("K45").Value = 4128.4287
$ws.Range("L45").Value = 3698.625
$ws.Range("M45").Value = -3751.4287
$ws.Range("N45").Value = -4452.625
$ws.Range("H64").Value = 504999.5
$ws.Range("J64").Value = 999999
$ws.Range("L64").Value = 999999
$ws.Range("N64").Value = -1000495
$ws.Range("H67").Value = 504999.5
$ws.Range("J67").Value = 999999
$ws.Range("L67").Value = 999999
$ws.Range("N67").Value = -1001715
$ws.Range("H122").Value = 2428.1667
$ws.Range("I122").Value = 2213.8
$ws.Range("K122").Value = 6641.400000000001
$ws.Range("M122").Value = -4191.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3154.0908
$ws.Range("I86").Value = 3288.5
$ws.Range("J86").Value = 2992.8
$ws.Range("K86").Value = 3288.5
$ws.Range("L86").Value = 2992.8
$ws.Range("M86").Value = -2165.5
$ws.Range("N86").Value = -5238.8
$ws.Range("H89").Value = 3154.0908
$ws.Range("I89").Value = 3288.5
$ws.Range("J89").Value = 2992.8
$ws.Range("K89").Value = 16442.5
$ws.Range("L89").Value = 14964
$ws.Range("M89").Value = -10826.5
$ws.Range("N89").Value = -26196
$ws.Range("H94").Value = 1665.2972
$ws.Range("I94").Value = 1227.0322
$ws.Range("J94").Value = 3929.6667
$ws.Range("K94").Value = 1227.0322
$ws.Range("L94").Value = 3929.6667
$ws.Range("M94").Value = -776.0322000000001
$ws.Range("N94").Value = -4831.6667
$ws.Range("H99").Value = 7178.5557
$ws.Range("I99").Value = 8443.929
$ws.Range("K99").Value = 8443.929
$ws.Range("M99").Value = -6945.929
$ws.Range("H105").Value = 2197.3872
$ws.Range("I105").Value = 1967.5
$ws.Range("K105").Value = 1967.5
$ws.Range("M105").Value = -220.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6000111
$ws.Range("I4").Value = 6000111
$ws.Range("K4").Value = 18000333
$ws.Range("M4").Value = -18000221
$ws.Range("H40").Value = 32.5
$ws.Range("I40").Value = 30
$ws.Range("K40").Value = 120
$ws.Range("M40").Value = -51
$ws.Range("H60").Value = 7499.6665
$ws.Range("I60").Value = 10249.5
$ws.Range("K60").Value = 30748.5
$ws.Range("M60").Value = -30497.5
$ws.Range("H92").Value = 836.8929000000001
$ws.Range("I92").Value = 499.8889
$ws.Range("J92").Value = 1443.5
$ws.Range("K92").Value = 1499.6667
$ws.Range("L92").Value = 4330.5
$ws.Range("M92").Value = -251.6667
$ws.Range("N92").Value = -6826.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3013.5
$ws.Range("J80").Value = 3209.5
$ws.Range("L80").Value = 3209.5
$ws.Range("N80").Value = -5205.5
$ws.Range("H83").Value = 3013.5
$ws.Range("J83").Value = 3209.5
$ws.Range("L83").Value = 16047.5
$ws.Range("N83").Value = -26031.5
$ws.Range("H122").Value = 5125.4165
$ws.Range("I122").Value = 5572.143
$ws.Range("K122").Value = 16716.429
$ws.Range("M122").Value = -14266.429
$ws.Range("H126").Value = 3410.2307
$ws.Range("I126").Value = 3179.75
$ws.Range("J126").Value = 3512.6667
$ws.Range("K126").Value = 9539.25
$ws.Range("L126").Value = 10538.0001
$ws.Range("M126").Value = -7069.25
$ws.Range("N126").Value = -15478.0001
$ws.Range("H132").Value = 7590.35
$ws.Range("I132").Value = 7884.579
$ws.Range("K132").Value = 23653.737
$ws.Range("M132").Value = -21123.737

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2575
$ws.Range("I7").Value = 2371.4285
$ws.Range("J7").Value = 4000
$ws.Range("K7").Value = 2371.4285
$ws.Range("L7").Value = 4000
$ws.Range("M7").Value = -2259.4285
$ws.Range("N7").Value = -4224
$ws.Range("H22").Value = 2641.6875
$ws.Range("I22").Value = 1174.1
$ws.Range("J22").Value = 3308.7727
$ws.Range("K22").Value = 1174.1
$ws.Range("L22").Value = 3308.7727
$ws.Range("M22").Value = -879.0999999999999
$ws.Range("N22").Value = -3898.7727
$ws.Range("H27").Value = 2641.6875
$ws.Range("I27").Value = 1174.1
$ws.Range("J27").Value = 3308.7727
$ws.Range("K27").Value = 1174.1
$ws.Range("L27").Value = 3308.7727
$ws.Range("M27").Value = -1067.1
$ws.Range("N27").Value = -3522.7727
$ws.Range("H40").Value = 2243.5
$ws.Range("I40").Value = 2262.3076
$ws.Range("J40").Value = 1999
$ws.Range("K40").Value = 2262.3076
$ws.Range("L40").Value = 1999
$ws.Range("M40").Value = -2126.3076
$ws.Range("N40").Value = -2271
$ws.Range("H68").Value = 15542.25
$ws.Range("I68").Value = 9658.333000000001
$ws.Range("K68").Value = 9658.333000000001
$ws.Range("M68").Value = -8909.333000000001
$ws.Range("H71").Value = 15542.25
$ws.Range("I71").Value = 9658.333000000001
$ws.Range("K71").Value = 48291.665
$ws.Range("M71").Value = -44547.665
$ws.Range("H93").Value = 4558.364
$ws.Range("I93").Value = 1908.8334
$ws.Range("J93").Value = 7737.8
$ws.Range("K93").Value = 1908.8334
$ws.Range("L93").Value = 7737.8
$ws.Range("M93").Value = -660.8334
$ws.Range("N93").Value = -10233.8
$ws.Range("H122").Value = 6156.643
$ws.Range("I122").Value = 4396.6
$ws.Range("K122").Value = 13189.8
$ws.Range("M122").Value = -10739.8
$ws.Range("H126").Value = 2575
$ws.Range("I126").Value = 2371.4285
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 7114.2855
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -4644.2855
$ws.Range("N126").Value = -16940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5999.5
$ws.Range("J14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("N14").Value = -8336
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5348
$ws.Range("H81").Value = 66736828
$ws.Range("I81").Value = 3060.125
$ws.Range("K81").Value = 6120.25
$ws.Range("M81").Value = -5059.25
$ws.Range("H84").Value = 66736828
$ws.Range("I84").Value = 3060.125
$ws.Range("K84").Value = 30601.25
$ws.Range("M84").Value = -25297.25
$ws.Range("H96").Value = 1306.1538
$ws.Range("I96").Value = 1361.8182
$ws.Range("K96").Value = 1361.8182
$ws.Range("M96").Value = 11.18180000000007
$ws.Range("H122").Value = 38495.5
$ws.Range("I122").Value = 3138.9583
$ws.Range("K122").Value = 9416.874899999999
$ws.Range("M122").Value = -6966.874899999999
$ws.Range("H126").Value = 2798.8
$ws.Range("I126").Value = 2623.5
$ws.Range("K126").Value = 7870.5
$ws.Range("M126").Value = -5400.5
